$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(161).Insert()

$ws.Range("A161").Value = 5
$ws.Range("B161").Value = "Macroferia Regional de Talca"
$ws.Range("C161").Value = "Maule"
$ws.Range("D161").Value = 44719
$ws.Range("E161").Value = 7
$ws.Range("F161").Value = "Fruta"
$ws.Range("G161").Value = 100108
$ws.Range("H161").Value = "Tropicales y subtropicales"
$ws.Range("I161").Value = 100108005
$ws.Range("J161").Value = "Piña"
$ws.Range("K161").Value = "Caramelo"
$ws.Range("L161").Value = "Segunda"
$ws.Range("M161").Value = 540
$ws.Range("N161").Value = 16000
$ws.Range("O161").Value = 16000
$ws.Range("P161").Value = 16000
$ws.Range("Q161").Value = "$/caja 14 unidades"
$ws.Range("R161").Value = "Ecuador"
$ws.Range("S161").Value = 1143
$ws.Range("T161").Value = 14
